$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo"), shifting "Tipo" to column E
$ws.Columns.Item(4).Insert()

# New header cell D1 = "MAE" (the inserted column already carries the header
# formatting, i.e. bold font + border, same as the other header cells)
$ws.Range("D1").Value = "MAE"

# New data cell D2 = MAE value for this row
$ws.Range("D2").Value = 0.8583964064824784
